$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format first so values such as
# "0.4700" or "5.990" keep trailing zeros instead of becoming numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "28.388.83"
$ws.Cells.Item(2, 5).Value = "  +3.32%  "
$ws.Cells.Item(3, 4).Value = "1.867.87"
$ws.Cells.Item(3, 5).Value = "  +1.58%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "339.11"
$ws.Cells.Item(5, 5).Value = "  +2.10%  "
$ws.Cells.Item(6, 5).Value = "  -0.03%  "
$ws.Cells.Item(7, 4).Value = "0.4700"
$ws.Cells.Item(7, 5).Value = "  +1.91%  "
$ws.Cells.Item(8, 4).Value = "0.3958"
$ws.Cells.Item(8, 5).Value = "  +3.63%  "
$ws.Cells.Item(9, 4).Value = "47.42"
$ws.Cells.Item(9, 5).Value = "  +2.44%  "
$ws.Cells.Item(10, 4).Value = "0.08012"
$ws.Cells.Item(10, 5).Value = "  +1.75%  "
$ws.Cells.Item(11, 4).Value = "0.9991"
$ws.Cells.Item(11, 5).Value = "  +2.28%  "
$ws.Cells.Item(12, 4).Value = "21.85"
$ws.Cells.Item(12, 5).Value = "  +3.20%  "
$ws.Cells.Item(13, 4).Value = "1.869.76"
$ws.Cells.Item(13, 5).Value = "  +1.37%  "
$ws.Cells.Item(14, 4).Value = "5.990"
$ws.Cells.Item(14, 5).Value = "  +1.73%  "
$ws.Cells.Item(15, 4).Value = "7.224"
$ws.Cells.Item(15, 5).Value = "  +2.87%  "
$ws.Cells.Item(16, 4).Value = "91.21"
$ws.Cells.Item(16, 5).Value = "  +3.89%  "
$ws.Cells.Item(17, 4).Value = "1.003"
$ws.Cells.Item(17, 5).Value = "  -0.05%  "
$ws.Cells.Item(18, 4).Value = "0.00001041"
$ws.Cells.Item(18, 5).Value = "  +1.36%  "
$ws.Cells.Item(19, 4).Value = "0.06614"
$ws.Cells.Item(19, 5).Value = "  -0.41%  "
$ws.Cells.Item(20, 4).Value = "17.53"
$ws.Cells.Item(20, 5).Value = "  +3.39%  "
$ws.Cells.Item(21, 4).Value = "1.001"
$ws.Cells.Item(21, 5).Value = "  -0.12%  "
$ws.Cells.Item(22, 4).Value = "28.401.82"
$ws.Cells.Item(22, 5).Value = "  +3.38%  "
$ws.Cells.Item(23, 4).Value = "5.453"
$ws.Cells.Item(23, 5).Value = "  +2.17%  "
$ws.Cells.Item(24, 5).Value = "  +1.72%  "
$ws.Cells.Item(25, 4).Value = "2.269"
$ws.Cells.Item(25, 5).Value = "  -0.83%  "
$ws.Cells.Item(26, 4).Value = "2.086.35"
$ws.Cells.Item(26, 5).Value = "  +1.10%  "
$ws.Cells.Item(27, 4).Value = "160.46"
$ws.Cells.Item(27, 5).Value = "  +2.20%  "
$ws.Cells.Item(28, 4).Value = "19.77"
$ws.Cells.Item(28, 5).Value = "  +2.22%  "
$ws.Cells.Item(29, 4).Value = "2.128"
$ws.Cells.Item(29, 5).Value = "  +2.85%  "
$ws.Cells.Item(30, 4).Value = "5.498"
$ws.Cells.Item(30, 5).Value = "  +3.37%  "
$ws.Cells.Item(31, 4).Value = "119.93"
$ws.Cells.Item(31, 5).Value = "  +1.08%  "
$ws.Cells.Item(32, 4).Value = "0.9658"
$ws.Cells.Item(32, 5).Value = "  +1.20%  "
$ws.Cells.Item(33, 4).Value = "0.09486"
$ws.Cells.Item(33, 5).Value = "  +2.12%  "
$ws.Cells.Item(34, 4).Value = "3.573"
$ws.Cells.Item(34, 5).Value = "  +0.41%  "
$ws.Cells.Item(35, 4).Value = "5.344"
$ws.Cells.Item(35, 5).Value = "  +2.20%  "
$ws.Cells.Item(36, 4).Value = "1.371"
$ws.Cells.Item(36, 5).Value = "  +3.72%  "
$ws.Cells.Item(37, 5).Value = "  +2.57%  "
$ws.Cells.Item(38, 4).Value = "0.02242"
$ws.Cells.Item(38, 5).Value = "  +2.30%  "
$ws.Cells.Item(39, 4).Value = "8.354"
$ws.Cells.Item(39, 5).Value = "  +3.47%  "
$ws.Cells.Item(40, 4).Value = "1.184"
$ws.Cells.Item(40, 5).Value = "  +1.83%  "
$ws.Cells.Item(41, 4).Value = "0.5939"
$ws.Cells.Item(41, 5).Value = "  +2.11%  "
$ws.Cells.Item(42, 4).Value = "1.001"
$ws.Cells.Item(42, 5).Value = "  -0.10%  "
$ws.Cells.Item(43, 5).Value = "  +1.47%  "
$ws.Cells.Item(44, 4).Value = "10.34"
$ws.Cells.Item(44, 5).Value = "  +2.52%  "
$ws.Cells.Item(45, 5).Value = "  +3.40%  "
$ws.Cells.Item(46, 4).Value = "0.5578"
$ws.Cells.Item(46, 5).Value = "  +1.64%  "
$ws.Cells.Item(47, 4).Value = "12.09"
$ws.Cells.Item(47, 5).Value = "  +1.41%  "
$ws.Cells.Item(48, 5).Value = "  +4.62%  "
$ws.Cells.Item(49, 4).Value = "0.06855"
$ws.Cells.Item(49, 5).Value = "  +2.98%  "
$ws.Cells.Item(50, 4).Value = "2.055"
$ws.Cells.Item(50, 5).Value = "  +15.79%  "
$ws.Cells.Item(51, 4).Value = "111.38"
$ws.Cells.Item(51, 5).Value = "  +1.25%  "
